# Auto-generated cell updates per commit: "chore: update Sheets via scheduled runner"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 133.78
$ws.Range("I15").Value = 133.78
$ws.Range("K15").Value = 401.34
$ws.Range("M15").Value = -232.34
# Row 93
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
# Row 129
$ws.Range("H129").Value = 3595.8108
$ws.Range("I129").Value = 17208.834
$ws.Range("J129").Value = 961.0323
$ws.Range("K129").Value = 51626.50199999999
$ws.Range("L129").Value = 2883.0969
$ws.Range("M129").Value = -46626.50199999999
$ws.Range("N129").Value = -12883.0969
# Row 137
$ws.Range("H137").Value = 1575
$ws.Range("I137").Value = 1500.0714
$ws.Range("K137").Value = 4500.2142
$ws.Range("M137").Value = -1950.2142
# Row 138
$ws.Range("H138").Value = 3684.1162
$ws.Range("I138").Value = 2752.9412
$ws.Range("J138").Value = 3913.5361
$ws.Range("K138").Value = 8258.8236
$ws.Range("L138").Value = 11740.6083
$ws.Range("M138").Value = -3118.8236
$ws.Range("N138").Value = -22020.6083

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
# Row 6
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
# Row 32
$ws.Range("H32").Value = 42577.375
$ws.Range("I32").Value = 18902.36
$ws.Range("J32").Value = 121000.875
$ws.Range("K32").Value = 18902.36
$ws.Range("L32").Value = 121000.875
$ws.Range("M32").Value = -18615.36
$ws.Range("N32").Value = -121574.875
# Row 61
$ws.Range("H61").Value = 2395.1177
$ws.Range("I61").Value = 2368.4666
$ws.Range("J61").Value = 2595
$ws.Range("K61").Value = 2368.4666
$ws.Range("L61").Value = 2595
$ws.Range("M61").Value = -2156.4666
$ws.Range("N61").Value = -3019
# Row 103
$ws.Range("H103").Value = 30000
$ws.Range("J103").Value = 30000
$ws.Range("L103").Value = 30000
$ws.Range("N103").Value = -32344
# Row 110
$ws.Range("H110").Value = 40085416
$ws.Range("I110").Value = 47720492
$ws.Range("J110").Value = 1269
$ws.Range("K110").Value = 47720492
$ws.Range("L110").Value = 1269
$ws.Range("M110").Value = -47718447
$ws.Range("N110").Value = -5359
# Row 132
$ws.Range("H132").Value = 8486407
$ws.Range("I132").Value = 9271841
$ws.Range("J132").Value = 3719.4
$ws.Range("K132").Value = 27815523
$ws.Range("L132").Value = 11158.2
$ws.Range("M132").Value = -27812993
$ws.Range("N132").Value = -16218.2
# Row 133
$ws.Range("H133").Value = 44990
$ws.Range("J133").Value = 44990
$ws.Range("L133").Value = 44990
$ws.Range("N133").Value = -50050
# Row 134
$ws.Range("H134").Value = 43330
$ws.Range("J134").Value = 43330
$ws.Range("L134").Value = 43330
$ws.Range("N134").Value = -53470
# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
# Row 136
$ws.Range("H136").Value = 2395.1177
$ws.Range("I136").Value = 2368.4666
$ws.Range("J136").Value = 2595
$ws.Range("K136").Value = 7105.399800000001
$ws.Range("L136").Value = 7785
$ws.Range("M136").Value = -4555.399800000001
$ws.Range("N136").Value = -12885
# Row 137
$ws.Range("H137").Value = 58500
$ws.Range("J137").Value = 58500
$ws.Range("L137").Value = 58500
$ws.Range("N137").Value = -68700
# Row 138
$ws.Range("H138").Value = 93000
$ws.Range("J138").Value = 93000
$ws.Range("L138").Value = 93000
$ws.Range("N138").Value = -103280
# Row 139
$ws.Range("H139").Value = 65000
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 1201.9
$ws.Range("J80").Value = 1514.3889
$ws.Range("L80").Value = 1514.3889
$ws.Range("N80").Value = -3510.3889
# Row 83
$ws.Range("H83").Value = 1201.9
$ws.Range("J83").Value = 1514.3889
$ws.Range("L83").Value = 7571.9445
$ws.Range("N83").Value = -17555.9445
# Row 141
$ws.Range("H141").Value = 45000
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 3501
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 3501
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 3501
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -3727
# Row 31
$ws.Range("H31").Value = 106643.86
$ws.Range("I31").Value = 2600
$ws.Range("K31").Value = 2600
$ws.Range("M31").Value = -2305
# Row 34
$ws.Range("H34").Value = 106643.86
$ws.Range("I34").Value = 2600
$ws.Range("K34").Value = 2600
$ws.Range("M34").Value = -2398
# Row 58
$ws.Range("H58").Value = 1296
$ws.Range("I58").Value = 1034.6487
$ws.Range("J58").Value = 2263
$ws.Range("K58").Value = 1034.6487
$ws.Range("L58").Value = 2263
$ws.Range("M58").Value = -831.6487
$ws.Range("N58").Value = -2669
# Row 75
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
# Row 78
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
# Row 136
$ws.Range("H136").Value = 1296
$ws.Range("I136").Value = 1034.6487
$ws.Range("J136").Value = 2263
$ws.Range("K136").Value = 3103.9461
$ws.Range("L136").Value = 6789
$ws.Range("M136").Value = -553.9461000000001
$ws.Range("N136").Value = -11889
# Row 138
$ws.Range("H138").Value = 92333.336
$ws.Range("J138").Value = 92333.336
$ws.Range("L138").Value = 92333.336
$ws.Range("N138").Value = -102613.336
# Row 139
$ws.Range("H139").Value = 44990
$ws.Range("J139").Value = 44990
$ws.Range("L139").Value = 44990
$ws.Range("N139").Value = -55270
# Row 140
$ws.Range("H140").Value = 59995
$ws.Range("J140").Value = 59995
$ws.Range("L140").Value = 59995
$ws.Range("N140").Value = -70355
# Row 141
$ws.Range("H141").Value = 68494.78
$ws.Range("J141").Value = 49932.668
$ws.Range("L141").Value = 49932.668
$ws.Range("N141").Value = -60292.668

$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Range("H39").Value = 5150.75
$ws.Range("I39").Value = 301.5
$ws.Range("J39").Value = 10000
$ws.Range("K39").Value = 904.5
$ws.Range("L39").Value = 30000
$ws.Range("M39").Value = -610.5
$ws.Range("N39").Value = -30588
# Row 55
$ws.Range("H55").Value = 26357.5
$ws.Range("J55").Value = 3880
$ws.Range("L55").Value = 11640
$ws.Range("N55").Value = -11994
# Row 98
$ws.Range("H98").Value = 126536
$ws.Range("J98").Value = 126536
$ws.Range("L98").Value = 379608
$ws.Range("N98").Value = -382604
# Row 114
$ws.Range("H114").Value = 2277.625
$ws.Range("I114").Value = 1350
$ws.Range("J114").Value = 2586.8333
$ws.Range("K114").Value = 4050
$ws.Range("L114").Value = 7760.499899999999
$ws.Range("M114").Value = -796
$ws.Range("N114").Value = -14268.4999
# Row 117
$ws.Range("H117").Value = 8765.200000000001
$ws.Range("I117").Value = 498
$ws.Range("J117").Value = 9355.714
$ws.Range("K117").Value = 1494
$ws.Range("L117").Value = 28067.142
$ws.Range("M117").Value = 1948
$ws.Range("N117").Value = -34951.142
# Row 129
$ws.Range("H129").Value = 178382.89
$ws.Range("I129").Value = 8610.571
$ws.Range("J129").Value = 232401.36
$ws.Range("K129").Value = 25831.713
$ws.Range("L129").Value = 697204.08
$ws.Range("M129").Value = -20831.713
$ws.Range("N129").Value = -707204.08
# Row 131
$ws.Range("H131").Value = 704983.75
$ws.Range("J131").Value = 735601.9
$ws.Range("L131").Value = 2206805.7
$ws.Range("N131").Value = -2216885.7

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 318158.1
$ws.Range("I102").Value = 2389.1
$ws.Range("K102").Value = 2389.1
$ws.Range("M102").Value = -767.0999999999999

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 79458.38
$ws.Range("J40").Value = 2632.7273
$ws.Range("L40").Value = 2632.7273
$ws.Range("N40").Value = -2904.7273

